# Sync attendance_reports: fix "Recorded By" text order for specific sessions
# ("dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com") on the
# sessions dated 07/12/2025, 08/12/2025, 14/12/2025, 21/12/2025, 22/12/2025,
# 27/12/2025 and 28/12/2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$targetDates = @("07/12/2025","08/12/2025","14/12/2025","21/12/2025","22/12/2025","27/12/2025","28/12/2025")

# Determine last used row in column E (Date)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dateVal = $ws.Cells.Item($r, 5).Value2
    if ($targetDates -contains $dateVal) {
        $gCell = $ws.Cells.Item($r, 7)
        if ($gCell.Value2 -eq $oldValue) {
            $gCell.Value = $newValue
        }
    }
}

$wb.Save()
